# Auto-generated Excel COM-interop script to apply scraped-schedule refresh
# Updates three worksheets (LP1912, LP1912-215, 6203-6173) with the latest scrape pass.
$wb = $excel.ActiveWorkbook

### Sheet: LP1912
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:28:27"
$ws.Cells.Item(3, 1).Value = "Total filas: 170"
$ws.Cells.Item(45, 1).Value = "08:14:55"
$ws.Cells.Item(45, 3).Value = "15_ABASTO"
$ws.Cells.Item(45, 4).Value = 15
$ws.Cells.Item(46, 1).Value = "07:51:40"
$ws.Cells.Item(46, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(46, 4).Value = 38
$ws.Cells.Item(55, 1).Value = "08:14:55"
$ws.Cells.Item(55, 3).Value = "215B_EL PATO"
$ws.Cells.Item(55, 4).Value = 39
$ws.Cells.Item(56, 1).Value = "08:49:06"
$ws.Cells.Item(56, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(56, 4).Value = 4
$ws.Cells.Item(66, 1).Value = "08:49:06"
$ws.Cells.Item(66, 3).Value = "14_ABASTO"
$ws.Cells.Item(66, 4).Value = 29
$ws.Cells.Item(67, 1).Value = "08:57:42"
$ws.Cells.Item(67, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(67, 4).Value = 21
$ws.Cells.Item(71, 1).Value = "08:14:55"
$ws.Cells.Item(71, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(71, 4).Value = 77
$ws.Cells.Item(72, 1).Value = "08:49:06"
$ws.Cells.Item(72, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(72, 4).Value = 42
$ws.Cells.Item(103, 1).Value = "10:32:07"
$ws.Cells.Item(103, 3).Value = "14_ABASTO"
$ws.Cells.Item(103, 4).Value = 43
$ws.Cells.Item(104, 1).Value = "11:01:19"
$ws.Cells.Item(104, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(104, 4).Value = 14
$ws.Cells.Item(117, 3).Value = "17_ROMERO"
$ws.Cells.Item(118, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(122, 1).Value = "11:56:32"
$ws.Cells.Item(122, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(122, 4).Value = 14
$ws.Cells.Item(123, 1).Value = "11:38:09"
$ws.Cells.Item(123, 3).Value = "15_ABASTO"
$ws.Cells.Item(123, 4).Value = 32
$ws.Cells.Item(137, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(138, 3).Value = "14_ABASTO"
$ws.Cells.Item(148, 1).Value = "12:43:13"
$ws.Cells.Item(148, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(148, 4).Value = 42
$ws.Cells.Item(149, 1).Value = "12:58:23"
$ws.Cells.Item(149, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(149, 4).Value = 27
$ws.Cells.Item(150, 1).Value = "13:28:27"
$ws.Cells.Item(150, 4).Value = 4
$ws.Cells.Item(151, 1).Value = "12:43:13"
$ws.Cells.Item(151, 3).Value = "14_ABASTO"
$ws.Cells.Item(151, 4).Value = 50
$ws.Cells.Item(152, 1).Value = "13:28:27"
$ws.Cells.Item(152, 3).Value = "215A_EL PATO"
$ws.Cells.Item(152, 4).Value = 5
$ws.Cells.Item(153, 1).Value = "13:28:27"
$ws.Cells.Item(153, 4).Value = 19
$ws.Cells.Item(155, 1).Value = "13:28:27"
$ws.Cells.Item(155, 4).Value = 26
$ws.Cells.Item(156, 1).Value = "13:28:27"
$ws.Cells.Item(156, 4).Value = 34
$ws.Cells.Item(159, 1).Value = "13:28:27"
$ws.Cells.Item(159, 2).Value = "14:06"
$ws.Cells.Item(159, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(159, 4).Value = 38
$ws.Cells.Item(160, 1).Value = "13:28:27"
$ws.Cells.Item(160, 2).Value = "14:14"
$ws.Cells.Item(160, 3).Value = "15_ABASTO"
$ws.Cells.Item(160, 4).Value = 46
$ws.Cells.Item(161, 1).Value = "13:28:27"
$ws.Cells.Item(161, 2).Value = "14:16"
$ws.Cells.Item(161, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(161, 4).Value = 48
$ws.Cells.Item(162, 1).Value = "13:28:27"
$ws.Cells.Item(162, 2).Value = "14:17"
$ws.Cells.Item(162, 4).Value = 49
$ws.Cells.Item(163, 2).Value = "14:17"
$ws.Cells.Item(163, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(163, 4).Value = 79
$ws.Cells.Item(164, 1).Value = "12:43:13"
$ws.Cells.Item(164, 2).Value = "14:18"
$ws.Cells.Item(164, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(164, 4).Value = 95
$ws.Cells.Item(165, 1).Value = "13:28:27"
$ws.Cells.Item(165, 2).Value = "14:27"
$ws.Cells.Item(165, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(165, 4).Value = 59
$ws.Cells.Item(166, 1).Value = "13:28:27"
$ws.Cells.Item(166, 2).Value = "14:32"
$ws.Cells.Item(166, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(166, 4).Value = 64
$ws.Cells.Item(167, 1).Value = "13:28:27"
$ws.Cells.Item(167, 2).Value = "14:34"
$ws.Cells.Item(167, 3).Value = "215C_EL PATO"
$ws.Cells.Item(167, 4).Value = 66
$ws.Cells.Item(168, 1).Value = "13:28:27"
$ws.Cells.Item(168, 2).Value = "14:39"
$ws.Cells.Item(168, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(168, 4).Value = 71
$ws.Cells.Item(169, 1).Value = "13:28:27"
$ws.Cells.Item(169, 2).Value = "14:47"
$ws.Cells.Item(169, 3).Value = "215B_EL PATO"
$ws.Cells.Item(169, 4).Value = 79
$ws.Cells.Item(169, 5).Value = "LP1912"
$ws.Cells.Item(170, 1).Value = "13:28:27"
$ws.Cells.Item(170, 2).Value = "14:51"
$ws.Cells.Item(170, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(170, 4).Value = 83
$ws.Cells.Item(170, 5).Value = "LP1912"
$ws.Cells.Item(171, 1).Value = "13:28:27"
$ws.Cells.Item(171, 2).Value = "14:51"
$ws.Cells.Item(171, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(171, 4).Value = 83
$ws.Cells.Item(171, 5).Value = "LP1912"
$ws.Cells.Item(172, 1).Value = "13:28:27"
$ws.Cells.Item(172, 2).Value = "14:54"
$ws.Cells.Item(172, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(172, 4).Value = 86
$ws.Cells.Item(172, 5).Value = "LP1912"
$ws.Cells.Item(173, 1).Value = "13:28:27"
$ws.Cells.Item(173, 2).Value = "15:02"
$ws.Cells.Item(173, 3).Value = "10_OLMOS"
$ws.Cells.Item(173, 4).Value = 94
$ws.Cells.Item(173, 5).Value = "LP1912"
$ws.Cells.Item(174, 1).Value = "13:28:27"
$ws.Cells.Item(174, 2).Value = "15:11"
$ws.Cells.Item(174, 3).Value = "14_ABASTO"
$ws.Cells.Item(174, 4).Value = 103
$ws.Cells.Item(174, 5).Value = "LP1912"
$ws.Cells.Item(175, 1).Value = "13:28:27"
$ws.Cells.Item(175, 2).Value = "15:13"
$ws.Cells.Item(175, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(175, 4).Value = 105
$ws.Cells.Item(175, 5).Value = "LP1912"

### Sheet: LP1912-215
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:28:27"
$ws.Cells.Item(30, 1).Value = "13:28:27"
$ws.Cells.Item(30, 4).Value = 5
$ws.Cells.Item(31, 1).Value = "13:28:27"
$ws.Cells.Item(31, 4).Value = 66
$ws.Cells.Item(32, 1).Value = "13:28:27"
$ws.Cells.Item(32, 4).Value = 79
$ws.Cells.Item(33, 1).Value = "13:28:27"
$ws.Cells.Item(33, 4).Value = 86

### Sheet: 6203-6173
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:28:27"
$ws.Cells.Item(3, 1).Value = "Total filas: 28"
$ws.Cells.Item(22, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(23, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(31, 1).Value = "13:28:27"
$ws.Cells.Item(31, 4).Value = 29
$ws.Cells.Item(32, 1).Value = "13:28:27"
$ws.Cells.Item(32, 2).Value = "14:27"
$ws.Cells.Item(32, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(32, 4).Value = 59
$ws.Cells.Item(32, 5).Value = "L6203"
$ws.Cells.Item(33, 1).Value = "13:28:27"
$ws.Cells.Item(33, 2).Value = "15:22"
$ws.Cells.Item(33, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(33, 4).Value = 114
$ws.Cells.Item(33, 5).Value = "L6173"

